$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.339.55"
$ws.Cells.Item(2, 5).Value = "  -0.33%  "
$ws.Cells.Item(3, 4).Value = "1.840.74"
$ws.Cells.Item(3, 5).Value = "  -0.84%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.19%  "
$ws.Cells.Item(5, 5).Value = "  -0.65%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.6290"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.76%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.14%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.07461"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.74%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.2901"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.86%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "24.82"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.85%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07737"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -0.28%  "
$ws.Cells.Item(12, 4).Value = "1.843.67"
$ws.Cells.Item(12, 5).Value = "  -0.68%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "4.979"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -1.13%  "
$ws.Cells.Item(14, 5).Value = "  -1.22%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001019"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.61%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "81.97"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.71%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "6.238"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.35%  "
$ws.Cells.Item(18, 4).Value = "29.332.14"
$ws.Cells.Item(18, 5).Value = "  -0.43%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "229.03"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.60%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "12.32"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.61%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.08%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "7.415"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.66%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.29%  "
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "158.79"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.23%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "8.467"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.21%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "0.1354"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -3.32%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "17.41"
$c.Style = "Normal"
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "0.06482"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +13.65%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.446"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +2.16%  "
$ws.Cells.Item(30, 5).Value = "  +0.31%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "4.064"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -2.32%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "4.064"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -0.02%  "
$ws.Cells.Item(33, 5).Value = "  +0.32%  "
$ws.Cells.Item(34, 5).Value = "  -1.75%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.6926"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.96%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "2.571"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.71%  "
$ws.Cells.Item(37, 5).Value = "  +1.44%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "2.814"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.38%  "
$ws.Cells.Item(39, 4).Value = "1.240.40"
$ws.Cells.Item(39, 5).Value = "  -0.84%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "6.752"
$c.Style = "Normal"
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.9300"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.01%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.19%  "
$ws.Cells.Item(43, 4).Value = "1.994.11"
$ws.Cells.Item(43, 5).Value = "  -1.16%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "100.81"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.70%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "65.60"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.76%  "
$ws.Cells.Item(46, 2).Value = "Aptos"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "7.048"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.67%  "
$ws.Cells.Item(47, 2).Value = "RenderToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "1.711"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.98%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "0.00000000116"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.64%  "
$ws.Cells.Item(49, 5).Value = "  -1.25%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "8.996"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -0.43%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.3901"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.84%  "